# ---------------------------------------------------------------------------
# HWP_hand_high.xlsx edit: refine the "Adjusted" sheet's values and add a new
# "Adjusted2" sheet (second attempt at the error-correction computation,
# using simpler cvxpy syntax) right after it.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Refine the numeric results already on the "Adjusted" sheet --------
$adjusted = $wb.Worksheets.Item("Adjusted")

$adjustedValues = @(
  @(0.03254268465275995, 0.9674573153472401, 0.1248031445349757, 0.12612628048467, 1.00000001037273),
  @(0.06242859853692478, 0.9375714014630752, 0.2081479833830336, 0.1233113590460544, 1.000000011269137),
  @(0.1209820810317784, 0.8790179189682216, 0.3009298829715466, 0.125644871166683, 1.000000022029525),
  @(0.2324891015093183, 0.7675108984906817, 0.3984428008581052, 0.1402899434485209, 1.000000029199465),
  @(0.4055817042816183, 0.5944182957183817, 0.466846017718665, 0.1521183653240594, 1.000000031790008),
  @(0.6274090824001115, 0.3725909175998885, 0.4597057432544561, 0.1497917590379689, 1.000000031469718),
  @(0.7918938665897139, 0.2081061334102861, 0.382557870606216, 0.1358214283533477, 1.000000028230799),
  @(0.8909518430700784, 0.1090481569299216, 0.2863550780381714, 0.1231155365457528, 1.000000019314174),
  @(0.9447900120694688, 0.05520998793053117, 0.1965202397394069, 0.1163685779702181, 1.000000010805606),
  @(0.9696170114597213, 0.03038298854027865, 0.1174042528649303, 0.1252042713255187, 1.000000011202573),
  @(0.9793136698290263, 0.02068633017097365, 0.03952013062523277, 0.1367354252113668, 1.000000022634675),
  @(0.968153358253151, 0.03184664174684898, -0.04500145749998, 0.1697271861563546, 1.000000031482759),
  @(0.9090446407686883, 0.09095535923131176, -0.166546406413095, 0.2344030552018472, 1.000000031837323),
  @(0.6997273631360486, 0.3002726368639514, -0.3064093236525103, 0.340913947421876, 1.000000025506467),
  @(0.2845527354297051, 0.7154472645702949, -0.3011473986993679, 0.3359951387941504, 1.000000025695038),
  @(0.079485205365186, 0.920514794634814, -0.1493813460766528, 0.2255050711686693, 1.000000032370438),
  @(0.03116296182507328, 0.9688370381749267, -0.04115648124722153, 0.1688134806348009, 1.000000031114654),
  @(0.02101318187701839, 0.9789868181229816, 0.04331263780899301, 0.1367320526911348, 1.000000021525365),
  @(0.03220308177013187, 0.9677969182298681, 0.1221860895279675, 0.1274229494984679, 1.000000010476775)
)

$r = 2
foreach ($rowVals in $adjustedValues) {
  $c = 2
  foreach ($val in $rowVals) {
    $adjusted.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

# --- 2. Add the new "Adjusted2" sheet right after "Adjusted" --------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $adjusted)
$newSheet.Name = "Adjusted2"

$headers = @("theta", "Jxx", "Jyy", "beta", "gamma", "trace")
for ($i = 0; $i -lt $headers.Length; $i++) {
  $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Reuse the exact header formatting (bold, centered, thin border) already
# used by "Adjusted" (and the other sheets) instead of re-building it from
# scratch, so the workbook keeps a single shared style definition.
$adjusted.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Match the page-margin convention used by the other sheets in the workbook
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$adjusted2Values = @(
  @(0, 0.00823926557676647, 1.004280151745791, 0.009748396696141807, 0, 1.008836571164085),
  @(5, 0.0400728232692661, 0.967590639033394, 0.1872324132046697, 0, 1.00794942901871),
  @(10, 0.1302413757149329, 0.8663714482967404, 0.3359884141175454, 0, 0.9933387312143561),
  @(15, 0.2669815085397383, 0.7272813522116881, 0.4384584699554157, 0, 0.9847089509283042),
  @(20, 0.4221613762681594, 0.5705130280537252, 0.4871778439718916, 0, 0.9783898461058586),
  @(25, 0.5982729450521631, 0.3950831970121937, 0.4862227700575049, 0, 0.9868464135875511),
  @(30, 0.7538380221324389, 0.2403338022287728, 0.425627151579058, 0, 0.9883490444288912),
  @(35, 0.8898654154169746, 0.1130531088849234, 0.3171034011236628, 0, 1.005750596992161),
  @(40, 0.971741668429856, 0.0286223414235253, 0.1666770057074869, 0, 1.000663557054632),
  @(45, 1.003100076847461, 0.007257711278065687, -0.006130973602136205, 0, 1.006337616218999),
  @(50, 0.9671013427791919, 0.04214208414597755, -0.1783546851674829, 0, 1.000681749903867),
  @(55, 0.8801643560147123, 0.1336077484880021, -0.3286418997007931, 0, 1.008551320532719),
  @(60, 0.7373475987130458, 0.2732739591968641, -0.4442503968562111, 0, 1.013076968316826),
  @(65, 0.585142931475698, 0.4275341595152552, -0.4981169113978479, 0, 1.021418622649452),
  @(70, 0.4200706042505858, 0.5917328120262253, -0.4958322201203499, 0, 1.018306214402867),
  @(75, 0.2559539526500945, 0.7562268978018446, -0.4337928423770912, 0, 1.0137440070314),
  @(80, 0.1281162565287108, 0.8875025188941184, -0.3244644612202032, 0, 1.014628869420169),
  @(85, 0.03779461798848963, 0.97363536423238, -0.1710884653446005, 0, 1.007936781580756),
  @(90, 0.009038262221567487, 1.005213329018905, 0.003154228577423192, 0, 1.010555425337093)
)

$r = 2
foreach ($rowVals in $adjusted2Values) {
  $c = 1
  foreach ($val in $rowVals) {
    $newSheet.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

# --- 3. Restore the originally active sheet/selection ---------------------
# Adding a sheet shifts UI focus to the new sheet; the source workbook keeps
# "measured" (the first tab) as the active one.
[void]$wb.Worksheets.Item("measured").Select()
[void]$wb.Worksheets.Item("measured").Range("A1").Select()
